{"js": "// Prepend \"Design: \" to every answer paragraph (style \"List Bullet\") in the\n// feedback table, matching the diff that turns e.g. \"Yes, 100% happy.\" into\n// \"Design: Yes, 100% happy.\" for all six Q&A rows.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst prefix = \"Design: \";\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.style === \"List Bullet\" && !paragraph.text.startsWith(prefix)) {\n    paragraph.insertText(prefix, \"Start\");\n  }\n}\nawait context.sync();\n", "ps1": "# Prepend \"Design: \" to every answer paragraph (style \"List Bullet\") in the\n# feedback table, matching the diff that turns e.g. \"Yes, 100% happy.\" into\n# \"Design: Yes, 100% happy.\" for all six Q&A rows.\n$d = $word.ActiveDocument\n\n$prefix = \"Design: \"\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"List Bullet\") {\n        $r = $p.Range\n        if ($r.Text.IndexOf($prefix) -ne 0) {\n            $r.InsertBefore($prefix)\n        }\n    }\n}\n"}
